$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-06 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-07 Monday", 2) | Out-Null
$d.Content.Find.Execute("665×2=", $true, $false, $false, $false, $false, $true, 1, $false, "263×6=", 2) | Out-Null
$d.Content.Find.Execute("169×9=", $true, $false, $false, $false, $false, $true, 1, $false, "203×4=", 2) | Out-Null
$d.Content.Find.Execute("353×3=", $true, $false, $false, $false, $false, $true, 1, $false, "453×4=", 2) | Out-Null
$d.Content.Find.Execute("634×5=", $true, $false, $false, $false, $false, $true, 1, $false, "549×2=", 2) | Out-Null
$d.Content.Find.Execute("240×8=", $true, $false, $false, $false, $false, $true, 1, $false, "486×6=", 2) | Out-Null
$d.Content.Find.Execute("518×5=", $true, $false, $false, $false, $false, $true, 1, $false, "977×8=", 2) | Out-Null
$d.Content.Find.Execute("212×6=", $true, $false, $false, $false, $false, $true, 1, $false, "759×8=", 2) | Out-Null
$d.Content.Find.Execute("669×6=", $true, $false, $false, $false, $false, $true, 1, $false, "855×7=", 2) | Out-Null
$d.Content.Find.Execute("786×3=", $true, $false, $false, $false, $false, $true, 1, $false, "729×9=", 2) | Out-Null
$d.Content.Find.Execute("727×2=", $true, $false, $false, $false, $false, $true, 1, $false, "839×2=", 2) | Out-Null
$d.Content.Find.Execute("708×9=", $true, $false, $false, $false, $false, $true, 1, $false, "354×2=", 2) | Out-Null
$d.Content.Find.Execute("178×7=", $true, $false, $false, $false, $false, $true, 1, $false, "483×4=", 2) | Out-Null
$d.Content.Find.Execute("701×4=", $true, $false, $false, $false, $false, $true, 1, $false, "350×8=", 2) | Out-Null
$d.Content.Find.Execute("546×3=", $true, $false, $false, $false, $false, $true, 1, $false, "525×4=", 2) | Out-Null
$d.Content.Find.Execute("354×7=", $true, $false, $false, $false, $false, $true, 1, $false, "455×6=", 2) | Out-Null
$d.Content.Find.Execute("139×5=", $true, $false, $false, $false, $false, $true, 1, $false, "743×4=", 2) | Out-Null
$d.Content.Find.Execute("530×4=", $true, $false, $false, $false, $false, $true, 1, $false, "882×4=", 2) | Out-Null
$d.Content.Find.Execute("283×2=", $true, $false, $false, $false, $false, $true, 1, $false, "635×8=", 2) | Out-Null
$d.Content.Find.Execute("573×2=", $true, $false, $false, $false, $false, $true, 1, $false, "292×4=", 2) | Out-Null
$d.Content.Find.Execute("210×2=", $true, $false, $false, $false, $false, $true, 1, $false, "859×3=", 2) | Out-Null
$d.Content.Find.Execute("868×5=", $true, $false, $false, $false, $false, $true, 1, $false, "845×9=", 2) | Out-Null
$d.Content.Find.Execute("900×9=", $true, $false, $false, $false, $false, $true, 1, $false, "791×7=", 2) | Out-Null
$d.Content.Find.Execute("796×6=", $true, $false, $false, $false, $false, $true, 1, $false, "311×3=", 2) | Out-Null
$d.Content.Find.Execute("715×2=", $true, $false, $false, $false, $false, $true, 1, $false, "809×7=", 2) | Out-Null
$d.Content.Find.Execute("123×7=", $true, $false, $false, $false, $false, $true, 1, $false, "903×3=", 2) | Out-Null

Write-Host "Replacements applied"
